$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.449.50"
$ws.Range("E2").Value = "  +3.23%  "
$ws.Range("D3").Value = "2.309.59"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'309.03"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'105.02"
$ws.Range("E6").Value = "  +8.55%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +7.85%  "
$ws.Range("D10").Value = "'36.22"
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("D11").Value = "'52.36"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'6.96"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "2.666.77"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "'15.11"
$ws.Range("E16").Value = "  +4.26%  "
$ws.Range("D17").Value = "2.310.90"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "'0.803"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").Value = "43.386.51"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "'11.95"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'6.16"
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").Value = "'67.87"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "'241.03"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "'2.61"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'24.87"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  +5.39%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'36.34"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'9.58"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'163.87"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'5.24"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'18.20"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("D37").Value = "'0.0736"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").Value = "'4.50"
$ws.Range("E39").Value = "  +9.40%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.106"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "'2.50"
$ws.Range("E43").Value = "  +13.45%  "
$ws.Range("D44").Value = "1.984.77"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'0.0290"
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("D46").Value = "'18.99"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("D47").Value = "'3.10"
$ws.Range("E47").Value = "  +6.55%  "
$ws.Range("D48").Value = "'10.21"
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("D49").Value = "'57.57"
$ws.Range("E49").Value = "  +6.81%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.58"
$ws.Range("E50").Value = "  +8.39%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.91"
$ws.Range("E51").Value = "  +0.84%  "
